$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sessions")

# Insert two new rows before row 44 (shifts old rows 44-48 down to 46-50)
$ws.Rows("44:45").Insert()
